$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '63.540.94'
$cell.ClearFormats()
$ws.Range("E2").Value = '  -0.05%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '3.068.22'
$cell.ClearFormats()
$ws.Range("E3").Value = '  -0.69%  '
$ws.Range("E4").Value = '  -0.30%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '593.52'
$cell.ClearFormats()
$ws.Range("E5").Value = '  +1.07%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '154.30'
$cell.ClearFormats()
$ws.Range("E6").Value = '  +1.31%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("E8").Value = '  -1.66%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '3.066.54'
$cell.ClearFormats()
$ws.Range("E9").Value = '  -0.35%  '
$ws.Range("E10").Value = '  -0.50%  '
$ws.Range("E11").Value = '  -0.03%  '
$ws.Range("E12").Value = '  -1.82%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '36.91'
$cell.ClearFormats()
$ws.Range("E13").Value = '  -1.61%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '0.0000238'
$cell.ClearFormats()
$ws.Range("E14").Value = '  -1.64%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.121'
$cell.ClearFormats()
$ws.Range("E15").Value = '  +1.32%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '3.572.50'
$cell.ClearFormats()
$ws.Range("E16").Value = '  -0.79%  '
$ws.Range("E17").Value = '  +0.82%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '63.492.01'
$cell.ClearFormats()
$ws.Range("E18").Value = '  -0.14%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '3.066.93'
$cell.ClearFormats()
$ws.Range("E19").Value = '  -0.62%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '490.60'
$cell.ClearFormats()
$ws.Range("E20").Value = '  +3.05%  '
$ws.Range("E21").Value = '  -1.44%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '0.708'
$cell.ClearFormats()
$ws.Range("E22").Value = '  -1.32%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '7.56'
$cell.ClearFormats()
$ws.Range("E23").Value = '  +0.10%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '2.47'
$cell.ClearFormats()
$ws.Range("E24").Value = '  +4.96%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '81.93'
$cell.ClearFormats()
$ws.Range("E25").Value = '  +0.15%  '
$ws.Range("E26").Value = '  -1.88%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '10.76'
$cell.ClearFormats()
$ws.Range("E27").Value = '  +11.08%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.ClearFormats()
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("E29").Value = '  +2.47%  '
$ws.Range("E30").Value = '  +0.49%  '
$cell = $ws.Range("B31")
$cell.NumberFormat = "@"
$cell.Value = 'ImmutableX'
$cell.ClearFormats()
$cell = $ws.Range("C31")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell.ClearFormats()
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '2.21'
$cell.ClearFormats()
$ws.Range("E31").Value = '  +0.85%  '
$cell = $ws.Range("B32")
$cell.NumberFormat = "@"
$cell.Value = 'FirstDigitalUSD'
$cell.ClearFormats()
$cell = $ws.Range("C32")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$cell.ClearFormats()
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.ClearFormats()
$ws.Range("E32").Value = '  -0.11%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '27.34'
$cell.ClearFormats()
$ws.Range("E33").Value = '  -0.34%  '
$ws.Range("E34").Value = '  -1.16%  '
$ws.Range("E35").Value = '  +1.08%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0823'
$cell.ClearFormats()
$ws.Range("E36").Value = '  -2.99%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '3.33'
$cell.ClearFormats()
$ws.Range("E37").Value = '  -0.93%  '
$ws.Range("E38").Value = '  -1.83%  '
$ws.Range("E39").Value = '  +0.56%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '9.29'
$cell.ClearFormats()
$ws.Range("E40").Value = '  -0.41%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '50.63'
$cell.ClearFormats()
$ws.Range("E41").Value = '  +0.07%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '438.05'
$cell.ClearFormats()
$ws.Range("E42").Value = '  -1.14%  '
$ws.Range("E43").Value = '  +2.60%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.113'
$cell.ClearFormats()
$ws.Range("E44").Value = '  +4.36%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.0364'
$cell.ClearFormats()
$ws.Range("E45").Value = '  -0.37%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '2.841.74'
$cell.ClearFormats()
$ws.Range("E46").Value = '  +0.78%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '39.23'
$cell.ClearFormats()
$ws.Range("E47").Value = '  +0.64%  '
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '130.04'
$cell.ClearFormats()
$ws.Range("E48").Value = '  -0.04%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '25.36'
$cell.ClearFormats()
$ws.Range("E49").Value = '  +0.93%  '
$ws.Range("E51").Value = '  -0.71%  '
